$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.982.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.63%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.351.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.569"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.91%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  -3.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.13%  "

$ws.Range("E11").Value = "  -0.15%  "

$ws.Range("E12").Value = "  -3.24%  "

$ws.Range("E13").Value = "  -0.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.714.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.355.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.97%  "

$ws.Range("E16").Value = "  -1.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.890.58"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.52%  "

$ws.Range("E20").Value = "  +2.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.71%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.10%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "39.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.79%  "

$ws.Range("E28").Value = "  -1.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.73%  "

$ws.Range("B31").Value = "LidoDAOToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +19.55%  "

$ws.Range("E32").Value = "  +5.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "146.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0773"
$ws.Range("D35").Style = "Normal"

$ws.Range("E36").Value = "  +0.53%  "

$ws.Range("E37").Value = "  +4.15%  "

$ws.Range("E38").Value = "  -2.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0300"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.877.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.73%  "

$ws.Range("E46").Value = "  -10.27%  "

$ws.Range("E47").Value = "  -6.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.584.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.61%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.85%  "
